# Update existing weekly data points (refreshed source numbers)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 11752
$ws.Range("D3").Value = 11497
$ws.Range("D7").Value = 11856
$ws.Range("D8").Value = 11143
$ws.Range("D9").Value = 713
$ws.Range("D11").Value = 564
$ws.Range("D12").Value = 11964
$ws.Range("D13").Value = 10641
$ws.Range("D14").Value = 1323
$ws.Range("D16").Value = 1044
$ws.Range("D17").Value = 12008
$ws.Range("D18").Value = 10090
$ws.Range("D19").Value = 1918
$ws.Range("D21").Value = 1502
$ws.Range("D22").Value = 12066
$ws.Range("D23").Value = 9662
$ws.Range("D24").Value = 2404
$ws.Range("D26").Value = 1884
$ws.Range("D27").Value = 12108
$ws.Range("D28").Value = 9203
$ws.Range("D29").Value = 2905
$ws.Range("D31").Value = 2296
$ws.Range("D32").Value = 12148
$ws.Range("D33").Value = 8795
$ws.Range("D34").Value = 3353
$ws.Range("D36").Value = 2647
$ws.Range("D37").Value = 12186
$ws.Range("D38").Value = 8271
$ws.Range("D39").Value = 3915
$ws.Range("D41").Value = 3114
$ws.Range("D42").Value = 12221
$ws.Range("D43").Value = 7704
$ws.Range("D44").Value = 4517
$ws.Range("D46").Value = 3615
$ws.Range("D47").Value = 12246
$ws.Range("D48").Value = 7156
$ws.Range("D49").Value = 5090
$ws.Range("D51").Value = 4075
$ws.Range("D52").Value = 12274
$ws.Range("D53").Value = 6514
$ws.Range("D54").Value = 5760
$ws.Range("D56").Value = 4579
$ws.Range("D57").Value = 12297
$ws.Range("D58").Value = 5903
$ws.Range("D59").Value = 6394
$ws.Range("D61").Value = 5106
$ws.Range("D62").Value = 12308
$ws.Range("D63").Value = 5737
$ws.Range("D64").Value = 6571
$ws.Range("D66").Value = 5254
$ws.Range("D67").Value = 12330
$ws.Range("D68").Value = 5553
$ws.Range("D69").Value = 6777
$ws.Range("D71").Value = 5412
$ws.Range("D72").Value = 12346
$ws.Range("D73").Value = 5117
$ws.Range("D74").Value = 7229
$ws.Range("D76").Value = 5813
$ws.Range("D77").Value = 12368
$ws.Range("D78").Value = 4639
$ws.Range("D79").Value = 7729
$ws.Range("D81").Value = 6257
$ws.Range("D82").Value = 12386
$ws.Range("D83").Value = 4120
$ws.Range("D84").Value = 8266
$ws.Range("D86").Value = 6758
$ws.Range("D87").Value = 12413
$ws.Range("D88").Value = 3714
$ws.Range("D89").Value = 8699
$ws.Range("D91").Value = 7161
$ws.Range("D92").Value = 12429
$ws.Range("D93").Value = 3376
$ws.Range("D94").Value = 9053
$ws.Range("D96").Value = 7519
$ws.Range("D97").Value = 12450
$ws.Range("D98").Value = 3065
$ws.Range("D99").Value = 9385
$ws.Range("D101").Value = 7836
$ws.Range("D102").Value = 12465
$ws.Range("D103").Value = 2789
$ws.Range("D104").Value = 9676
$ws.Range("D106").Value = 8127
$ws.Range("D107").Value = 12475
$ws.Range("D108").Value = 2502
$ws.Range("D109").Value = 9973
$ws.Range("D111").Value = 8432

# Append the new week (202510, last day of week 2025-03-09) for every
# Variable category, matching the existing pattern of 5 rows per week.
$ws.Range("B111:D111").Copy() | Out-Null
$ws.Range("B112:D116").PasteSpecial(-4122) | Out-Null

$ws.Range("A112").Value = 202510
$ws.Range("B112").Value = 45725
$ws.Range("C112").Value = "farms_total_count"
$ws.Range("D112").Value = 12496

$ws.Range("A113").Value = 202510
$ws.Range("B113").Value = 45725
$ws.Range("C113").Value = "farms_to_examine_count"
$ws.Range("D113").Value = 2265

$ws.Range("A114").Value = 202510
$ws.Range("B114").Value = 45725
$ws.Range("C114").Value = "farms_examined_count"
$ws.Range("D114").Value = 10231

$ws.Range("A115").Value = 202510
$ws.Range("B115").Value = 45725
$ws.Range("C115").Value = "farms_examined_positive_count"
$ws.Range("D115").Value = 1534

$ws.Range("A116").Value = 202510
$ws.Range("B116").Value = 45725
$ws.Range("C116").Value = "farms_examined_negative_count"
$ws.Range("D116").Value = 8697

# Update the sheet view to match where the author left the selection
$ws.Range("G15").Select() | Out-Null
